$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.353.61"
$ws.Range("E2").Value = "  -0.75%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.176.98"
$ws.Range("E3").Value = "  -1.71%  "

# Row 4
$ws.Range("E4").Value = "  -0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.48"
$ws.Range("E5").Value = "  +4.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.611"
$ws.Range("E6").Value = "  -0.06%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.03"
$ws.Range("E7").Value = "  -2.36%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.581"
$ws.Range("E9").Value = "  -3.81%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.76"
$ws.Range("E10").Value = "  -3.33%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0908"
$ws.Range("E11").Value = "  -2.09%  "

# Row 12
$ws.Range("E12").Value = "  +0.11%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.75"
$ws.Range("E13").Value = "  -2.00%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.501.99"
$ws.Range("E14").Value = "  -1.75%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.17"
$ws.Range("E15").Value = "  -3.56%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.172.75"
$ws.Range("E16").Value = "  -2.23%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.768"
$ws.Range("E17").Value = "  -4.17%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.260.81"
$ws.Range("E18").Value = "  -0.59%  "

# Row 19
$ws.Range("E19").Value = "  -2.73%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.57"
$ws.Range("E20").Value = "  -0.26%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.84"
$ws.Range("E21").Value = "  -1.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.24"
$ws.Range("E22").Value = "  -1.83%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.33"
$ws.Range("E23").Value = "  -4.86%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("E24").Value = "  -2.33%  "

# Row 25
$ws.Range("E25").Value = "  +0.08%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.43"
$ws.Range("E26").Value = "  -4.29%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.38"
$ws.Range("E27").Value = "  +0.31%  "

# Row 28
$ws.Range("E28").Value = "  -2.40%  "

# Row 29
$ws.Range("E29").Value = "  -2.20%  "

# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.76"
$ws.Range("E30").Value = "  +0.85%  "

# Row 31
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.78"
$ws.Range("E31").Value = "  -1.76%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.97"
$ws.Range("E32").Value = "  -1.49%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0813"
$ws.Range("E33").Value = "  +2.29%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.10"
$ws.Range("E34").Value = "  -3.40%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.120"
$ws.Range("E35").Value = "  -1.82%  "

# Row 36
$ws.Range("E36").Value = "  -1.67%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0335"
$ws.Range("E37").Value = "  +3.87%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.18"
$ws.Range("E38").Value = "  -4.91%  "

# Row 39
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.71"
$ws.Range("E39").Value = "  -5.44%  "

# Row 40
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.04"
$ws.Range("E40").Value = "  -4.13%  "

# Row 41
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "59.09"
$ws.Range("E41").Value = "  -1.91%  "

# Row 42
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.194"
$ws.Range("E42").Value = "  -1.26%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.13"
$ws.Range("E43").Value = "  -6.73%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.12"
$ws.Range("E44").Value = "  +1.98%  "

# Row 45
$ws.Range("E45").Value = "  +7.32%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0969"
$ws.Range("E46").Value = "  -2.03%  "

# Row 47
$ws.Range("B47").Value = "WOONetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.459"
$ws.Range("E47").Value = "  +6.72%  "

# Row 48
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.15"
$ws.Range("E48").Value = "  -4.84%  "

# Row 49
$ws.Range("E49").Value = "  -2.40%  "

# Row 50
$ws.Range("E50").Value = "  -1.38%  "

# Row 51
$ws.Range("E51").Value = "  +0.26%  "
